$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H76").Value = 3609.0908
$ws.Range("I76").Value = 3128.5715
$ws.Range("J76").Value = 4450
$ws.Range("K76").Value = 3128.5715
$ws.Range("L76").Value = 4450
$ws.Range("M76").Value = -2813.5715
$ws.Range("N76").Value = -5080
$ws.Range("H79").Value = 3609.0908
$ws.Range("I79").Value = 3128.5715
$ws.Range("J79").Value = 4450
$ws.Range("K79").Value = 3128.5715
$ws.Range("L79").Value = 4450
$ws.Range("M79").Value = -2036.5715
$ws.Range("N79").Value = -6634
$ws.Range("H135").Value = 663.6774
$ws.Range("I135").Value = 683.2174
$ws.Range("J135").Value = 607.5
$ws.Range("K135").Value = 6148.9566
$ws.Range("L135").Value = 5467.5
$ws.Range("M135").Value = -3613.9566
$ws.Range("N135").Value = -10537.5
$ws.Range("H137").Value = 1836.4032
$ws.Range("I137").Value = 1313.5385
$ws.Range("J137").Value = 2214.0278
$ws.Range("K137").Value = 3940.6155
$ws.Range("L137").Value = 6642.0834
$ws.Range("M137").Value = -1390.6155
$ws.Range("N137").Value = -11742.0834
$ws.Range("H138").Value = 5058.8706
$ws.Range("I138").Value = 4597.8
$ws.Range("J138").Value = 5105.9185
$ws.Range("K138").Value = 13793.4
$ws.Range("L138").Value = 15317.7555
$ws.Range("M138").Value = -8653.400000000001
$ws.Range("N138").Value = -25597.7555

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 24087.887
$ws.Range("I2").Value = 1042.303
$ws.Range("K2").Value = 1042.303
$ws.Range("M2").Value = -929.3030000000001
$ws.Range("H32").Value = 665.17
$ws.Range("I32").Value = 604.2526
$ws.Range("K32").Value = 604.2526
$ws.Range("M32").Value = -317.2526
$ws.Range("H61").Value = 1926.738
$ws.Range("I61").Value = 903.06665
$ws.Range("J61").Value = 2495.4443
$ws.Range("K61").Value = 903.06665
$ws.Range("L61").Value = 2495.4443
$ws.Range("M61").Value = -691.06665
$ws.Range("N61").Value = -2919.4443
$ws.Range("H63").Value = 2656.6667
$ws.Range("I63").Value = 2411.111
$ws.Range("J63").Value = 3393.3333
$ws.Range("K63").Value = 2411.111
$ws.Range("L63").Value = 3393.3333
$ws.Range("M63").Value = -1725.111
$ws.Range("N63").Value = -4765.3333
$ws.Range("H66").Value = 2656.6667
$ws.Range("I66").Value = 2411.111
$ws.Range("J66").Value = 3393.3333
$ws.Range("K66").Value = 12055.555
$ws.Range("L66").Value = 16966.6665
$ws.Range("M66").Value = -8623.555
$ws.Range("N66").Value = -23830.6665
$ws.Range("H88").Value = 2034.3
$ws.Range("I88").Value = 1599.8
$ws.Range("K88").Value = 1599.8
$ws.Range("M88").Value = -1193.8
$ws.Range("H91").Value = 2034.3
$ws.Range("I91").Value = 1599.8
$ws.Range("K91").Value = 1599.8
$ws.Range("M91").Value = -195.8
$ws.Range("H103").Value = 46362
$ws.Range("J103").Value = 46362
$ws.Range("L103").Value = 46362
$ws.Range("N103").Value = -48706
$ws.Range("H116").Value = 24087.887
$ws.Range("I116").Value = 1042.303
$ws.Range("K116").Value = 1042.303
$ws.Range("M116").Value = 1251.697
$ws.Range("H136").Value = 1926.738
$ws.Range("I136").Value = 903.06665
$ws.Range("J136").Value = 2495.4443
$ws.Range("K136").Value = 2709.19995
$ws.Range("L136").Value = 7486.3329
$ws.Range("M136").Value = -159.1999500000002
$ws.Range("N136").Value = -12586.3329

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 24087.887
$ws.Range("I3").Value = 1042.303
$ws.Range("K3").Value = 1042.303
$ws.Range("M3").Value = -928.3030000000001
$ws.Range("H82").Value = 22254
$ws.Range("I82").Value = 4599.75
$ws.Range("J82").Value = 36377.4
$ws.Range("K82").Value = 4599.75
$ws.Range("L82").Value = 36377.4
$ws.Range("M82").Value = -4216.75
$ws.Range("N82").Value = -37143.4
$ws.Range("H85").Value = 22254
$ws.Range("I85").Value = 4599.75
$ws.Range("J85").Value = 36377.4
$ws.Range("K85").Value = 4599.75
$ws.Range("L85").Value = 36377.4
$ws.Range("M85").Value = -3273.75
$ws.Range("N85").Value = -39029.4
$ws.Range("H107").Value = 19611276
$ws.Range("I107").Value = 23812288
$ws.Range("J107").Value = 6559.3335
$ws.Range("K107").Value = 23812288
$ws.Range("L107").Value = 6559.3335
$ws.Range("M107").Value = -23810368
$ws.Range("N107").Value = -10399.3335

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 18088.746
$ws.Range("I31").Value = 41185.6
$ws.Range("J31").Value = 2893.4473
$ws.Range("K31").Value = 41185.6
$ws.Range("L31").Value = 2893.4473
$ws.Range("M31").Value = -40890.6
$ws.Range("N31").Value = -3483.4473
$ws.Range("H34").Value = 18088.746
$ws.Range("I34").Value = 41185.6
$ws.Range("J34").Value = 2893.4473
$ws.Range("K34").Value = 41185.6
$ws.Range("L34").Value = 2893.4473
$ws.Range("M34").Value = -40983.6
$ws.Range("N34").Value = -3297.4473
$ws.Range("H86").Value = 4107.4443
$ws.Range("I86").Value = 3877
$ws.Range("J86").Value = 4136.25
$ws.Range("K86").Value = 3877
$ws.Range("L86").Value = 4136.25
$ws.Range("M86").Value = -2754
$ws.Range("N86").Value = -6382.25
$ws.Range("H89").Value = 4107.4443
$ws.Range("I89").Value = 3877
$ws.Range("J89").Value = 4136.25
$ws.Range("K89").Value = 19385
$ws.Range("L89").Value = 20681.25
$ws.Range("M89").Value = -13769
$ws.Range("N89").Value = -31913.25
$ws.Range("H96").Value = 21333.334
$ws.Range("J96").Value = 21333.334
$ws.Range("L96").Value = 21333.334
$ws.Range("N96").Value = -26825.334

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H68").Value = 1945.904
$ws.Range("I68").Value = 1286.75
$ws.Range("J68").Value = 2356.0444
$ws.Range("K68").Value = 3860.25
$ws.Range("L68").Value = 7068.1332
$ws.Range("M68").Value = -3049.25
$ws.Range("N68").Value = -8690.1332
$ws.Range("H71").Value = 1945.904
$ws.Range("I71").Value = 1286.75
$ws.Range("J71").Value = 2356.0444
$ws.Range("K71").Value = 11580.75
$ws.Range("L71").Value = 21204.3996
$ws.Range("M71").Value = -7524.75
$ws.Range("N71").Value = -29316.3996
$ws.Range("H86").Value = 1230.6666
$ws.Range("I86").Value = 1746
$ws.Range("J86").Value = 200
$ws.Range("K86").Value = 5238
$ws.Range("L86").Value = 600
$ws.Range("M86").Value = -4052
$ws.Range("N86").Value = -2972
$ws.Range("H89").Value = 1230.6666
$ws.Range("I89").Value = 1746
$ws.Range("J89").Value = 200
$ws.Range("K89").Value = 15714
$ws.Range("L89").Value = 1800
$ws.Range("M89").Value = -9786
$ws.Range("N89").Value = -13656
$ws.Range("H98").Value = 86014.08
$ws.Range("I98").Value = 850
$ws.Range("J98").Value = 101498.45
$ws.Range("K98").Value = 2550
$ws.Range("L98").Value = 304495.35
$ws.Range("M98").Value = -1052
$ws.Range("N98").Value = -307491.35

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 1697.8182
$ws.Range("I102").Value = 1806.9615
$ws.Range("K102").Value = 1806.9615
$ws.Range("M102").Value = -184.9614999999999
$ws.Range("H132").Value = 1896.138
$ws.Range("I132").Value = 1378.619
$ws.Range("J132").Value = 3254.625
$ws.Range("K132").Value = 4135.857
$ws.Range("L132").Value = 9763.875
$ws.Range("M132").Value = -1605.857
$ws.Range("N132").Value = -14823.875

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 127637.5
$ws.Range("I16").Value = 334700.34
$ws.Range("J16").Value = 3399.8
$ws.Range("K16").Value = 334700.34
$ws.Range("L16").Value = 3399.8
$ws.Range("M16").Value = -334530.34
$ws.Range("N16").Value = -3739.8
$ws.Range("H61").Value = 2162.125
$ws.Range("I61").Value = 2042.4286
$ws.Range("K61").Value = 2042.4286
$ws.Range("M61").Value = -1840.4286
$ws.Range("H68").Value = 2673.8333
$ws.Range("I68").Value = 1620.5454
$ws.Range("K68").Value = 1620.5454
$ws.Range("M68").Value = -871.5454
$ws.Range("H71").Value = 2673.8333
$ws.Range("I71").Value = 1620.5454
$ws.Range("K71").Value = 8102.727
$ws.Range("M71").Value = -4358.727
$ws.Range("H113").Value = 2162.125
$ws.Range("I113").Value = 2042.4286
$ws.Range("K113").Value = 2042.4286
$ws.Range("M113").Value = 127.5714
$ws.Range("H132").Value = 9911.888999999999
$ws.Range("I132").Value = 16251
$ws.Range("J132").Value = 4840.6
$ws.Range("K132").Value = 48753
$ws.Range("L132").Value = 14521.8
$ws.Range("M132").Value = -46223
$ws.Range("N132").Value = -19581.8

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H107").Value = 39973.84
$ws.Range("I107").Value = 10921
$ws.Range("J107").Value = 92797.17999999999
$ws.Range("K107").Value = 32763
$ws.Range("L107").Value = 278391.54
$ws.Range("M107").Value = -30843
$ws.Range("N107").Value = -282231.54
$ws.Range("H122").Value = 1983.1
$ws.Range("I122").Value = 1918
$ws.Range("J122").Value = 2135
$ws.Range("K122").Value = 5754
$ws.Range("L122").Value = 6405
$ws.Range("M122").Value = -3304
$ws.Range("N122").Value = -11305
